$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58: 2026-01-23 (serial 46045), Error Count 47
$ws.Range("A57").Copy($ws.Range("A58"))
$ws.Range("A58").Value = 46045
$ws.Range("B58").Value = 47

# Row 59: 2026-01-26 (serial 46048), Error Count 66
$ws.Range("A57").Copy($ws.Range("A59"))
$ws.Range("A59").Value = 46048
$ws.Range("B59").Value = 66

# Update view: select A58:B59 (new bottom rows) so the active cell / selection
# matches the post-edit state; the sheet's topLeftCell scroll position is
# window-chrome state that isn't exposed through this COM surface.
$ws.Range("A58:B59").Select() | Out-Null
